$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.472738
$ws.Range("H2").Value = 1.418214
$ws.Range("I2").Value = 0.0327564895931267
$ws.Range("J2").Value = 0.03397138804734427
$ws.Range("M2").Value = 63.46725166666666
$ws.Range("N2").Value = 190.401755
$ws.Range("O2").Value = 0.2354497988808272
$ws.Range("P2").Value = 0.2397164477183668
$ws.Range("Q2").Value = 30.00338161839667
$ws.Range("R2").Value = 270.03043456557
$ws.Range("S2").Value = 0.007712508886743592
$ws.Range("T2").Value = 0.008143500466771553

$ws.Range("G3").Value = 0.472738
$ws.Range("H3").Value = 1.418214
$ws.Range("I3").Value = 0.0327564895931267
$ws.Range("J3").Value = 0.03397138804734427
$ws.Range("O3").Value = 0.1779985000094065
$ws.Range("P3").Value = 0.1812240584798697
$ws.Range("Q3").Value = 22.68235924885
$ws.Range("R3").Value = 204.14123323965
$ws.Range("S3").Value = 0.005830606013150285
$ws.Range("T3").Value = 0.006156432814134265

$ws.Range("G4").Value = 0.472738
$ws.Range("H4").Value = 1.418214
$ws.Range("I4").Value = 0.0327564895931267
$ws.Range("J4").Value = 0.03397138804734427
$ws.Range("M4").Value = 64.53809233333334
$ws.Range("N4").Value = 193.614277
$ws.Range("O4").Value = 0.2394223865221556
$ws.Range("P4").Value = 0.243761023683841
$ws.Range("Q4").Value = 30.50960869347534
$ws.Range("R4").Value = 274.586478241278
$ws.Range("S4").Value = 0.007842636912474549
$ws.Range("T4").Value = 0.008280900326381641

$ws.Range("G5").Value = 0.472738
$ws.Range("H5").Value = 1.418214
$ws.Range("I5").Value = 0.0327564895931267
$ws.Range("J5").Value = 0.03397138804734427
$ws.Range("M5").Value = 14.3933435
$ws.Range("N5").Value = 28.786687
$ws.Range("O5").Value = 0.0533961963580272
$ws.Range("P5").Value = 0.03624253541791403
$ws.Range("Q5").Value = 6.804280419503001
$ws.Range("R5").Value = 40.825682517018
$ws.Range("S5").Value = 0.001749071950314268
$ws.Range("T5").Value = 0.001231209234501576

$ws.Range("G6").Value = 0.472738
$ws.Range("H6").Value = 1.418214
$ws.Range("I6").Value = 0.0327564895931267
$ws.Range("J6").Value = 0.03397138804734427
$ws.Range("M6").Value = 79.17795566666666
$ws.Range("N6").Value = 237.533867
$ws.Range("O6").Value = 0.2937331182295834
$ws.Range("P6").Value = 0.2990559347000084
$ws.Range("Q6").Value = 37.43042840594867
$ws.Range("R6").Value = 336.873855653538
$ws.Range("S6").Value = 0.009621665830444001
$ws.Range("T6").Value = 0.01015934520555523

$ws.Range("I7").Value = 0.822180234441485
$ws.Range("J7").Value = 0.8526739017519405
$ws.Range("M7").Value = 63.46725166666666
$ws.Range("N7").Value = 190.401755
$ws.Range("O7").Value = 0.2354497988808272
$ws.Range("P7").Value = 0.2397164477183668
$ws.Range("Q7").Value = 753.0778676060222
$ws.Range("R7").Value = 6777.700808454199
$ws.Range("S7").Value = 0.193582170843039
$ws.Range("T7").Value = 0.2043999587901348

$ws.Range("I8").Value = 0.822180234441485
$ws.Range("J8").Value = 0.8526739017519405
$ws.Range("O8").Value = 0.1779985000094065
$ws.Range("P8").Value = 0.1812240584798697
$ws.Range("S8").Value = 0.1463468484679665
$ws.Range("T8").Value = 0.1545250250353523

$ws.Range("I9").Value = 0.822180234441485
$ws.Range("J9").Value = 0.8526739017519405
$ws.Range("M9").Value = 64.53809233333334
$ws.Range("N9").Value = 193.614277
$ws.Range("O9").Value = 0.2394223865221556
$ws.Range("P9").Value = 0.243761023683841
$ws.Range("Q9").Value = 765.7840488982979
$ws.Range("R9").Value = 6892.056440084681
$ws.Range("S9").Value = 0.1968483538813258
$ws.Range("T9").Value = 0.2078486631595479

$ws.Range("I10").Value = 0.822180234441485
$ws.Range("J10").Value = 0.8526739017519405
$ws.Range("M10").Value = 14.3933435
$ws.Range("N10").Value = 28.786687
$ws.Range("O10").Value = 0.0533961963580272
$ws.Range("P10").Value = 0.03624253541791403
$ws.Range("Q10").Value = 170.7858485448467
$ws.Range("R10").Value = 1024.71509126908
$ws.Range("S10").Value = 0.04390129723992637
$ws.Range("T10").Value = 0.03090306408417565

$ws.Range("I11").Value = 0.822180234441485
$ws.Range("J11").Value = 0.8526739017519405
$ws.Range("M11").Value = 79.17795566666666
$ws.Range("N11").Value = 237.533867
$ws.Range("O11").Value = 0.2937331182295834
$ws.Range("P11").Value = 0.2990559347000084
$ws.Range("Q11").Value = 939.4950064644755
$ws.Range("R11").Value = 8455.45505818028
$ws.Range("S11").Value = 0.2415015640092273
$ws.Range("T11").Value = 0.2549971906827296

$ws.Range("G12").Value = 0.37892
$ws.Range("H12").Value = 1.13676
$ws.Range("I12").Value = 0.02625574638939025
$ws.Range("J12").Value = 0.02722954016579943
$ws.Range("M12").Value = 63.46725166666666
$ws.Range("N12").Value = 190.401755
$ws.Range("O12").Value = 0.2354497988808272
$ws.Range("P12").Value = 0.2397164477183668
$ws.Range("Q12").Value = 24.04901100153333
$ws.Range("R12").Value = 216.4410990138
$ws.Range("S12").Value = 0.00618191020684794
$ws.Range("T12").Value = 0.006527368641550027

$ws.Range("G13").Value = 0.37892
$ws.Range("H13").Value = 1.13676
$ws.Range("I13").Value = 0.02625574638939025
$ws.Range("J13").Value = 0.02722954016579943
$ws.Range("O13").Value = 0.1779985000094065
$ws.Range("P13").Value = 0.1812240584798697
$ws.Range("Q13").Value = 18.180894209
$ws.Range("R13").Value = 163.628047881
$ws.Range("S13").Value = 0.004673483473938853
$ws.Range("T13").Value = 0.004934647779386797

$ws.Range("G14").Value = 0.37892
$ws.Range("H14").Value = 1.13676
$ws.Range("I14").Value = 0.02625574638939025
$ws.Range("J14").Value = 0.02722954016579943
$ws.Range("M14").Value = 64.53809233333334
$ws.Range("N14").Value = 193.614277
$ws.Range("O14").Value = 0.2394223865221556
$ws.Range("P14").Value = 0.243761023683841
$ws.Range("Q14").Value = 24.45477394694667
$ws.Range("R14").Value = 220.09296552252
$ws.Range("S14").Value = 0.006286213460468285
$ws.Range("T14").Value = 0.006637500585255534

$ws.Range("G15").Value = 0.37892
$ws.Range("H15").Value = 1.13676
$ws.Range("I15").Value = 0.02625574638939025
$ws.Range("J15").Value = 0.02722954016579943
$ws.Range("M15").Value = 14.3933435
$ws.Range("N15").Value = 28.786687
$ws.Range("O15").Value = 0.0533961963580272
$ws.Range("P15").Value = 0.03624253541791403
$ws.Range("Q15").Value = 5.45392571902
$ws.Range("R15").Value = 32.72355431412
$ws.Range("S15").Value = 0.001401956989734445
$ws.Range("T15").Value = 0.0009868675738724986

$ws.Range("G16").Value = 0.37892
$ws.Range("H16").Value = 1.13676
$ws.Range("I16").Value = 0.02625574638939025
$ws.Range("J16").Value = 0.02722954016579943
$ws.Range("M16").Value = 79.17795566666666
$ws.Range("N16").Value = 237.533867
$ws.Range("O16").Value = 0.2937331182295834
$ws.Range("P16").Value = 0.2990559347000084
$ws.Range("Q16").Value = 30.00211096121333
$ws.Range("R16").Value = 270.01899865092
$ws.Range("S16").Value = 0.007712182258400722
$ws.Range("T16").Value = 0.00814315558573457

$ws.Range("G17").Value = 1.548357
$ws.Range("H17").Value = 3.096714
$ws.Range("I17").Value = 0.1072872076222874
$ws.Range("J17").Value = 0.0741775733180209
$ws.Range("M17").Value = 63.46725166666666
$ws.Range("N17").Value = 190.401755
$ws.Range("O17").Value = 0.2354497988808272
$ws.Range("P17").Value = 0.2397164477183668
$ws.Range("Q17").Value = 98.269963388845
$ws.Range("R17").Value = 589.61978033307
$ws.Range("S17").Value = 0.02526075145715311
$ws.Range("T17").Value = 0.01778158437616468

$ws.Range("G18").Value = 1.548357
$ws.Range("H18").Value = 3.096714
$ws.Range("I18").Value = 0.1072872076222874
$ws.Range("J18").Value = 0.0741775733180209
$ws.Range("O18").Value = 0.1779985000094065
$ws.Range("P18").Value = 0.1812240584798697
$ws.Range("Q18").Value = 74.29144625452501
$ws.Range("R18").Value = 445.7486775271501
$ws.Range("S18").Value = 0.01909696202696491
$ws.Range("T18").Value = 0.01344276088487984

$ws.Range("G19").Value = 1.548357
$ws.Range("H19").Value = 3.096714
$ws.Range("I19").Value = 0.1072872076222874
$ws.Range("J19").Value = 0.0741775733180209
$ws.Range("M19").Value = 64.53809233333334
$ws.Range("N19").Value = 193.614277
$ws.Range("O19").Value = 0.2394223865221556
$ws.Range("P19").Value = 0.243761023683841
$ws.Range("Q19").Value = 99.92800703096302
$ws.Range("R19").Value = 599.5680421857782
$ws.Range("S19").Value = 0.02568695929222605
$ws.Range("T19").Value = 0.01808160120638394

$ws.Range("G20").Value = 1.548357
$ws.Range("H20").Value = 3.096714
$ws.Range("I20").Value = 0.1072872076222874
$ws.Range("J20").Value = 0.0741775733180209
$ws.Range("M20").Value = 14.3933435
$ws.Range("N20").Value = 28.786687
$ws.Range("O20").Value = 0.0533961963580272
$ws.Range("P20").Value = 0.03624253541791403
$ws.Range("Q20").Value = 22.2860341616295
$ws.Range("R20").Value = 89.144136646518
$ws.Range("S20").Value = 0.005728728804904088
$ws.Range("T20").Value = 0.002688383328193287

$ws.Range("G21").Value = 1.548357
$ws.Range("H21").Value = 3.096714
$ws.Range("I21").Value = 0.1072872076222874
$ws.Range("J21").Value = 0.0741775733180209
$ws.Range("M21").Value = 79.17795566666666
$ws.Range("N21").Value = 237.533867
$ws.Range("O21").Value = 0.2937331182295834
$ws.Range("P21").Value = 0.2990559347000084
$ws.Range("Q21").Value = 122.595741902173
$ws.Range("R21").Value = 735.5744514130381
$ws.Range("S21").Value = 0.03151380604103919
$ws.Range("T21").Value = 0.02218324352239914

$ws.Range("G22").Value = 0.16626
$ws.Range("H22").Value = 0.49878
$ws.Range("I22").Value = 0.01152032195371061
$ws.Range("J22").Value = 0.01194759671689489
$ws.Range("M22").Value = 63.46725166666666
$ws.Range("N22").Value = 190.401755
$ws.Range("O22").Value = 0.2354497988808272
$ws.Range("P22").Value = 0.2397164477183668
$ws.Range("Q22").Value = 10.5520652621
$ws.Range("R22").Value = 94.9685873589
$ws.Range("S22").Value = 0.002712457487043541
$ws.Range("T22").Value = 0.002864035443745665

$ws.Range("G23").Value = 0.16626
$ws.Range("H23").Value = 0.49878
$ws.Range("I23").Value = 0.01152032195371061
$ws.Range("J23").Value = 0.01194759671689489
$ws.Range("O23").Value = 0.1779985000094065
$ws.Range("P23").Value = 0.1812240584798697
$ws.Range("Q23").Value = 7.9772919645
$ws.Range("R23").Value = 71.7956276805
$ws.Range("S23").Value = 0.002050600027385923
$ws.Range("T23").Value = 0.002165191966116459

$ws.Range("G24").Value = 0.16626
$ws.Range("H24").Value = 0.49878
$ws.Range("I24").Value = 0.01152032195371061
$ws.Range("J24").Value = 0.01194759671689489
$ws.Range("M24").Value = 64.53809233333334
$ws.Range("N24").Value = 193.614277
$ws.Range("O24").Value = 0.2394223865221556
$ws.Range("P24").Value = 0.243761023683841
$ws.Range("Q24").Value = 10.73010323134
$ws.Range("R24").Value = 96.57092908206
$ws.Range("S24").Value = 0.002758222975660976
$ws.Range("T24").Value = 0.002912358406271997

$ws.Range("G25").Value = 0.16626
$ws.Range("H25").Value = 0.49878
$ws.Range("I25").Value = 0.01152032195371061
$ws.Range("J25").Value = 0.01194759671689489
$ws.Range("M25").Value = 14.3933435
$ws.Range("N25").Value = 28.786687
$ws.Range("O25").Value = 0.0533961963580272
$ws.Range("P25").Value = 0.03624253541791403
$ws.Range("Q25").Value = 2.39303729031
$ws.Range("R25").Value = 14.35822374186
$ws.Range("S25").Value = 0.000615141373148023
$ws.Range("T25").Value = 0.0004330111971710166

$ws.Range("G26").Value = 0.16626
$ws.Range("H26").Value = 0.49878
$ws.Range("I26").Value = 0.01152032195371061
$ws.Range("J26").Value = 0.01194759671689489
$ws.Range("M26").Value = 79.17795566666666
$ws.Range("N26").Value = 237.533867
$ws.Range("O26").Value = 0.2937331182295834
$ws.Range("P26").Value = 0.2990559347000084
$ws.Range("Q26").Value = 13.16412690914
$ws.Range("R26").Value = 118.47714218226
$ws.Range("S26").Value = 0.003383900090472142
$ws.Range("T26").Value = 0.003572999703589754
